# Updated Master data as per 16th May Refresh
# Adds three new user_detail_h rows (Nikola Tesla, Graham Bell, Albert Miles),
# fixes the is_active (column I) left-alignment style that was missing on the
# last existing row, and resets the sheet's out-of-view selection anchor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix pre-existing style glitch on row 33 -------------------------------
# Every other data row has the boolean is_active cell (column I) left
# aligned; row 33 was missing that formatting. Re-apply it so it matches the
# rest of the column.
$ws.Range("I33").HorizontalAlignment = -4131  # xlLeft

# --- New data rows (34-36) --------------------------------------------------
$newRows = @(
    @{ Row = 34; Id = 110033; Uin = 9317596771; Name = "Nikola Tesla"; Email = "nikola.tesla@xyz.com"; Mobile = 818876434 },
    @{ Row = 35; Id = 110034; Uin = 9317596772; Name = "Graham Bell";  Email = "graham.bell@xyz.com";  Mobile = 818876435 },
    @{ Row = 36; Id = 110035; Uin = 9317596773; Name = "Albert Miles"; Email = "albert.miles@xyz.com"; Mobile = 818876436 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Id        # A: id
    $ws.Cells.Item($row, 2).Value = $r.Uin       # B: uin
    $ws.Cells.Item($row, 3).Value = $r.Name      # C: name
    $ws.Cells.Item($row, 4).Value = $r.Email     # D: email
    $ws.Cells.Item($row, 5).Value = $r.Mobile    # E: mobile
    $ws.Cells.Item($row, 6).Value = "ACT"        # F: status_code
    $ws.Cells.Item($row, 7).Value = "eng"        # G: lang_code
    $ws.Cells.Item($row, 8).Value = "PWD"        # H: last_login_method
    $ws.Cells.Item($row, 9).Value = $true        # I: is_active
    $ws.Cells.Item($row, 10).Value = "superadmin" # J: cr_by
    $ws.Cells.Item($row, 11).Value = "now()"     # K: cr_dtimes
    $ws.Cells.Item($row, 12).Value = "now()"     # L: eff_dtimes

    # Match the existing column formatting used throughout the table.
    $ws.Cells.Item($row, 9).HorizontalAlignment = -4131  # xlLeft, like column I elsewhere
}

# --- Reset the selection anchor (matches the saved view) -------------------
$null = $ws.Range("M1:XFD1048576").Select()
